$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DeviceList")
$ws.Columns.Item(5).Delete()

foreach ($cf in $ws.Range("B2:G2").FormatConditions) {
    $cf.ModifyAppliesToRange($ws.Range("B2:G2"))
}

$ws.Range("B14").Select() | Out-Null
